$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 5700.4
$ws.Cells.Item(40, 9).Value = 6750
$ws.Cells.Item(40, 10).Value = 1502
$ws.Cells.Item(40, 11).Value = 6750
$ws.Cells.Item(40, 12).Value = 1502
$ws.Cells.Item(40, 13).Value = -6575
$ws.Cells.Item(40, 14).Value = -1852

$ws.Cells.Item(80, 8).Value = 6979438
$ws.Cells.Item(80, 9).Value = 13644.444
$ws.Cells.Item(80, 10).Value = 13945232
$ws.Cells.Item(80, 11).Value = 40933.33199999999
$ws.Cells.Item(80, 12).Value = 41835696
$ws.Cells.Item(80, 13).Value = -39935.33199999999
$ws.Cells.Item(80, 14).Value = -41837692

$ws.Cells.Item(83, 8).Value = 6979438
$ws.Cells.Item(83, 9).Value = 13644.444
$ws.Cells.Item(83, 10).Value = 13945232
$ws.Cells.Item(83, 11).Value = 122799.996
$ws.Cells.Item(83, 12).Value = 125507088
$ws.Cells.Item(83, 13).Value = -117807.996
$ws.Cells.Item(83, 14).Value = -125517072

$ws.Cells.Item(86, 8).Value = 675677.4
$ws.Cells.Item(86, 9).Value = 1979.8
$ws.Cells.Item(86, 10).Value = 1012526.2
$ws.Cells.Item(86, 11).Value = 1979.8
$ws.Cells.Item(86, 12).Value = 1012526.2
$ws.Cells.Item(86, 13).Value = -856.8
$ws.Cells.Item(86, 14).Value = -1014772.2

$ws.Cells.Item(89, 8).Value = 675677.4
$ws.Cells.Item(89, 9).Value = 1979.8
$ws.Cells.Item(89, 10).Value = 1012526.2
$ws.Cells.Item(89, 11).Value = 9899
$ws.Cells.Item(89, 12).Value = 5062631
$ws.Cells.Item(89, 13).Value = -4283
$ws.Cells.Item(89, 14).Value = -5073863

$ws.Cells.Item(137, 8).Value = 13898676
$ws.Cells.Item(137, 9).Value = 2976898.5
$ws.Cells.Item(137, 10).Value = 52124900
$ws.Cells.Item(137, 11).Value = 8930695.5
$ws.Cells.Item(137, 12).Value = 156374700
$ws.Cells.Item(137, 13).Value = -8928145.5
$ws.Cells.Item(137, 14).Value = -156379800

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 14).ClearContents()

$ws.Cells.Item(61, 8).Value = 2138983
$ws.Cells.Item(61, 9).Value = 1069295.6
$ws.Cells.Item(61, 11).Value = 1069295.6
$ws.Cells.Item(61, 13).Value = -1069083.6

$ws.Cells.Item(74, 8).Value = 85979976
$ws.Cells.Item(74, 9).Value = 78022540
$ws.Cells.Item(74, 10).Value = 106669300
$ws.Cells.Item(74, 11).Value = 78022540
$ws.Cells.Item(74, 12).Value = 106669300
$ws.Cells.Item(74, 13).Value = -78021666
$ws.Cells.Item(74, 14).Value = -106671048

$ws.Cells.Item(77, 8).Value = 85979976
$ws.Cells.Item(77, 9).Value = 78022540
$ws.Cells.Item(77, 10).Value = 106669300
$ws.Cells.Item(77, 11).Value = 390112700
$ws.Cells.Item(77, 12).Value = 533346500
$ws.Cells.Item(77, 13).Value = -390108332
$ws.Cells.Item(77, 14).Value = -533355236

$ws.Cells.Item(136, 8).Value = 2138983
$ws.Cells.Item(136, 9).Value = 1069295.6
$ws.Cells.Item(136, 11).Value = 3207886.8
$ws.Cells.Item(136, 13).Value = -3205336.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(7, 8).Value = 6680.4
$ws.Cells.Item(7, 9).Value = 7600.75
$ws.Cells.Item(7, 10).Value = 2999
$ws.Cells.Item(7, 11).Value = 7600.75
$ws.Cells.Item(7, 12).Value = 2999
$ws.Cells.Item(7, 13).Value = -7487.75
$ws.Cells.Item(7, 14).Value = -3225

$ws.Cells.Item(94, 8).Value = 1576.3334
$ws.Cells.Item(94, 9).Value = 1156.7142
$ws.Cells.Item(94, 10).Value = 2310.6667
$ws.Cells.Item(94, 11).Value = 1156.7142
$ws.Cells.Item(94, 12).Value = 2310.6667
$ws.Cells.Item(94, 13).Value = -705.7141999999999
$ws.Cells.Item(94, 14).Value = -3212.6667

$ws.Cells.Item(105, 8).Value = 1908.091
$ws.Cells.Item(105, 9).Value = 1927
$ws.Cells.Item(105, 10).Value = 1875
$ws.Cells.Item(105, 11).Value = 1927
$ws.Cells.Item(105, 12).Value = 1875
$ws.Cells.Item(105, 13).Value = -180
$ws.Cells.Item(105, 14).Value = -5369

$ws.Cells.Item(134, 8).Value = 10715768
$ws.Cells.Item(134, 9).Value = 11629097
$ws.Cells.Item(134, 10).Value = 5105314
$ws.Cells.Item(134, 11).Value = 34887291
$ws.Cells.Item(134, 12).Value = 15315942
$ws.Cells.Item(134, 13).Value = -34884756
$ws.Cells.Item(134, 14).Value = -15321012

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(2, 8).Value = 100
$ws.Cells.Item(2, 9).Value = 100
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 100
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = 13
$ws.Cells.Item(2, 14).ClearContents()

$ws.Cells.Item(31, 8).Value = 2370568.2
$ws.Cells.Item(31, 9).Value = 1303550.6
$ws.Cells.Item(31, 10).Value = 5215948
$ws.Cells.Item(31, 11).Value = 1303550.6
$ws.Cells.Item(31, 12).Value = 5215948
$ws.Cells.Item(31, 13).Value = -1303255.6
$ws.Cells.Item(31, 14).Value = -5216538

$ws.Cells.Item(34, 8).Value = 2370568.2
$ws.Cells.Item(34, 9).Value = 1303550.6
$ws.Cells.Item(34, 10).Value = 5215948
$ws.Cells.Item(34, 11).Value = 1303550.6
$ws.Cells.Item(34, 12).Value = 5215948
$ws.Cells.Item(34, 13).Value = -1303348.6
$ws.Cells.Item(34, 14).Value = -5216352

$ws.Cells.Item(58, 8).Value = 2901952.2
$ws.Cells.Item(58, 9).Value = 1988195
$ws.Cells.Item(58, 11).Value = 1988195
$ws.Cells.Item(58, 13).Value = -1987992

$ws.Cells.Item(62, 8).Value = 2624.9167
$ws.Cells.Item(62, 9).Value = 2589.9
$ws.Cells.Item(62, 10).Value = 2800
$ws.Cells.Item(62, 11).Value = 2589.9
$ws.Cells.Item(62, 12).Value = 2800
$ws.Cells.Item(62, 13).Value = -1965.9
$ws.Cells.Item(62, 14).Value = -4048

$ws.Cells.Item(65, 8).Value = 2624.9167
$ws.Cells.Item(65, 9).Value = 2589.9
$ws.Cells.Item(65, 10).Value = 2800
$ws.Cells.Item(65, 11).Value = 12949.5
$ws.Cells.Item(65, 12).Value = 14000
$ws.Cells.Item(65, 13).Value = -9829.5
$ws.Cells.Item(65, 14).Value = -20240

$ws.Cells.Item(132, 8).Value = 1615247.2
$ws.Cells.Item(132, 9).Value = 2779496.5
$ws.Cells.Item(132, 10).Value = 3209.5386
$ws.Cells.Item(132, 11).Value = 8338489.5
$ws.Cells.Item(132, 12).Value = 9628.6158
$ws.Cells.Item(132, 13).Value = -8335959.5
$ws.Cells.Item(132, 14).Value = -14688.6158

$ws.Cells.Item(134, 8).Value = 2676778.5
$ws.Cells.Item(134, 9).Value = 11483.6
$ws.Cells.Item(134, 10).Value = 8007368.5
$ws.Cells.Item(134, 11).Value = 34450.8
$ws.Cells.Item(134, 12).Value = 24022105.5
$ws.Cells.Item(134, 13).Value = -31915.8
$ws.Cells.Item(134, 14).Value = -24027175.5

$ws.Cells.Item(136, 8).Value = 2901952.2
$ws.Cells.Item(136, 9).Value = 1988195
$ws.Cells.Item(136, 11).Value = 5964585
$ws.Cells.Item(136, 13).Value = -5962035

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(33, 8).Value = 87.8
$ws.Cells.Item(33, 9).Value = 52.5
$ws.Cells.Item(33, 10).Value = 111.333336
$ws.Cells.Item(33, 11).Value = 315
$ws.Cells.Item(33, 12).Value = 668.000016
$ws.Cells.Item(33, 13).Value = -32
$ws.Cells.Item(33, 14).Value = -1234.000016

$ws.Cells.Item(62, 8).Value = 200100740
$ws.Cells.Item(62, 9).Value = 1225
$ws.Cells.Item(62, 10).Value = 333500400
$ws.Cells.Item(62, 11).Value = 3675
$ws.Cells.Item(62, 12).Value = 1000501200
$ws.Cells.Item(62, 13).Value = -2989
$ws.Cells.Item(62, 14).Value = -1000502572

$ws.Cells.Item(63, 8).Value = 3300
$ws.Cells.Item(63, 9).Value = 2166.6667
$ws.Cells.Item(63, 10).Value = 5000
$ws.Cells.Item(63, 11).Value = 6500.000100000001
$ws.Cells.Item(63, 12).Value = 15000
$ws.Cells.Item(63, 13).Value = -5751.000100000001
$ws.Cells.Item(63, 14).Value = -16498

$ws.Cells.Item(65, 8).Value = 200100740
$ws.Cells.Item(65, 9).Value = 1225
$ws.Cells.Item(65, 10).Value = 333500400
$ws.Cells.Item(65, 11).Value = 11025
$ws.Cells.Item(65, 12).Value = 3001503600
$ws.Cells.Item(65, 13).Value = -7593
$ws.Cells.Item(65, 14).Value = -3001510464

$ws.Cells.Item(66, 8).Value = 3300
$ws.Cells.Item(66, 9).Value = 2166.6667
$ws.Cells.Item(66, 10).Value = 5000
$ws.Cells.Item(66, 11).Value = 19500.0003
$ws.Cells.Item(66, 12).Value = 45000
$ws.Cells.Item(66, 13).Value = -15756.0003
$ws.Cells.Item(66, 14).Value = -52488

$ws.Cells.Item(87, 8).Value = 1500
$ws.Cells.Item(87, 9).Value = 1500
$ws.Cells.Item(87, 11).Value = 4500
$ws.Cells.Item(87, 13).Value = -3252

$ws.Cells.Item(88, 8).Value = 5214.0386
$ws.Cells.Item(88, 9).Value = 1005
$ws.Cells.Item(88, 10).Value = 5382.4
$ws.Cells.Item(88, 11).Value = 3015
$ws.Cells.Item(88, 12).Value = 16147.2
$ws.Cells.Item(88, 13).Value = -2587
$ws.Cells.Item(88, 14).Value = -17003.2

$ws.Cells.Item(90, 8).Value = 1500
$ws.Cells.Item(90, 9).Value = 1500
$ws.Cells.Item(90, 11).Value = 13500
$ws.Cells.Item(90, 13).Value = -7260

$ws.Cells.Item(91, 8).Value = 5214.0386
$ws.Cells.Item(91, 9).Value = 1005
$ws.Cells.Item(91, 10).Value = 5382.4
$ws.Cells.Item(91, 11).Value = 3015
$ws.Cells.Item(91, 12).Value = 16147.2
$ws.Cells.Item(91, 13).Value = -1533
$ws.Cells.Item(91, 14).Value = -19111.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 28213352
$ws.Cells.Item(132, 9).Value = 49524410
$ws.Cells.Item(132, 10).Value = 12991169
$ws.Cells.Item(132, 11).Value = 148573230
$ws.Cells.Item(132, 12).Value = 38973507
$ws.Cells.Item(132, 13).Value = -148570700
$ws.Cells.Item(132, 14).Value = -38978567

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(10, 8).Value = 660
$ws.Cells.Item(10, 9).Value = 690
$ws.Cells.Item(10, 10).Value = 600
$ws.Cells.Item(10, 11).Value = 690
$ws.Cells.Item(10, 12).Value = 600
$ws.Cells.Item(10, 13).Value = -550
$ws.Cells.Item(10, 14).Value = -880

$ws.Cells.Item(82, 8).Value = 2988.5806
$ws.Cells.Item(82, 9).Value = 1423
$ws.Cells.Item(82, 10).Value = 4119.278
$ws.Cells.Item(82, 11).Value = 1423
$ws.Cells.Item(82, 12).Value = 4119.278
$ws.Cells.Item(82, 13).Value = -1062
$ws.Cells.Item(82, 14).Value = -4841.278

$ws.Cells.Item(85, 8).Value = 2988.5806
$ws.Cells.Item(85, 9).Value = 1423
$ws.Cells.Item(85, 10).Value = 4119.278
$ws.Cells.Item(85, 11).Value = 1423
$ws.Cells.Item(85, 12).Value = 4119.278
$ws.Cells.Item(85, 13).Value = -175
$ws.Cells.Item(85, 14).Value = -6615.278

$ws.Cells.Item(132, 8).Value = 1758892.4
$ws.Cells.Item(132, 9).Value = 2226930.8
$ws.Cells.Item(132, 10).Value = 3748.25
$ws.Cells.Item(132, 11).Value = 6680792.399999999
$ws.Cells.Item(132, 12).Value = 11244.75
$ws.Cells.Item(132, 13).Value = -6678262.399999999
$ws.Cells.Item(132, 14).Value = -16304.75

$ws.Cells.Item(136, 8).Value = 3923576.2
$ws.Cells.Item(136, 9).Value = 6537990
$ws.Cells.Item(136, 11).Value = 19613970
$ws.Cells.Item(136, 13).Value = -19611420

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 14777.895
$ws.Cells.Item(81, 9).Value = 729
$ws.Cells.Item(81, 10).Value = 22973.084
$ws.Cells.Item(81, 11).Value = 1458
$ws.Cells.Item(81, 12).Value = 45946.168
$ws.Cells.Item(81, 13).Value = -397
$ws.Cells.Item(81, 14).Value = -48068.168

$ws.Cells.Item(84, 8).Value = 14777.895
$ws.Cells.Item(84, 9).Value = 729
$ws.Cells.Item(84, 10).Value = 22973.084
$ws.Cells.Item(84, 11).Value = 7290
$ws.Cells.Item(84, 12).Value = 229730.84
$ws.Cells.Item(84, 13).Value = -1986
$ws.Cells.Item(84, 14).Value = -240338.84

$ws.Cells.Item(100, 8).Value = 7311
$ws.Cells.Item(100, 9).Value = 9243.091
$ws.Cells.Item(100, 10).Value = 226.66667
$ws.Cells.Item(100, 11).Value = 18486.182
$ws.Cells.Item(100, 12).Value = 453.33334
$ws.Cells.Item(100, 13).Value = -17945.182
$ws.Cells.Item(100, 14).Value = -1535.33334

$ws.Cells.Item(132, 8).Value = 1622379.1
$ws.Cells.Item(132, 9).Value = 1279969.9
$ws.Cells.Item(132, 10).Value = 2333536.5
$ws.Cells.Item(132, 11).Value = 3839909.7
$ws.Cells.Item(132, 12).Value = 7000609.5
$ws.Cells.Item(132, 13).Value = -3837379.7
$ws.Cells.Item(132, 14).Value = -7005669.5

$ws.Cells.Item(136, 8).Value = 14224.5
$ws.Cells.Item(136, 9).Value = 10578.2
$ws.Cells.Item(136, 10).Value = 20301.666
$ws.Cells.Item(136, 11).Value = 31734.6
$ws.Cells.Item(136, 12).Value = 60904.99800000001
$ws.Cells.Item(136, 13).Value = -29184.6
$ws.Cells.Item(136, 14).Value = -66004.99800000001
